$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 84), pushing the
# existing rows 84-102 down to 86-104. Excel's Rows.Insert() carries the
# formatting (incl. the date-style on column D) from the row being pushed
# down, matching the s="2" style already used throughout column D.
$ws.Rows(84).Insert()
$ws.Rows(84).Insert()

# Fill the two new rows (84-85) with the new weekly price records.
$ws.Range("A84").Value = 1
$ws.Range("B84").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C84").Value = "Arica y Parinacota"
$ws.Range("D84").Value = 44642
$ws.Range("E84").Value = 15
$ws.Range("F84").Value = 100112008
$ws.Range("G84").Value = "Coliflor"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 650
$ws.Range("L84").Value = 750
$ws.Range("M84").Value = 700
$ws.Range("N84").Value = "$/unidad"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 700
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"

$ws.Range("A85").Value = 1
$ws.Range("B85").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C85").Value = "Arica y Parinacota"
$ws.Range("D85").Value = 44642
$ws.Range("E85").Value = 15
$ws.Range("F85").Value = 100112008
$ws.Range("G85").Value = "Coliflor"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Tercera"
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 500
$ws.Range("L85").Value = 600
$ws.Range("M85").Value = 550
$ws.Range("N85").Value = "$/unidad"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 550
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"
